# Updates 'F' column ('想去人数' / interested-count) values across the
# '展览', '演出', and '全部类型' sheets to match the refreshed scrape,
# per commit 'Update gh-pages to output generated at 456a3b4'.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 7779
$ws.Range("F3").Value = 7779
$ws.Range("F5").Value = 7934
$ws.Range("F9").Value = 6794
$ws.Range("F10").Value = 3413
$ws.Range("F12").Value = 3750
$ws.Range("F13").Value = 50
$ws.Range("F14").Value = 56
$ws.Range("F15").Value = 51
$ws.Range("F16").Value = 77
$ws.Range("F17").Value = 85
$ws.Range("F18").Value = 481
$ws.Range("F20").Value = 63
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = 339
$ws.Range("F25").Value = 3910
$ws.Range("F26").Value = 120
$ws.Range("F28").Value = 1013
$ws.Range("F29").Value = 342
$ws.Range("F30").Value = 1530
$ws.Range("F31").Value = 85
$ws.Range("F32").Value = 69
$ws.Range("F33").Value = 2809
$ws.Range("F34").Value = 1971
$ws.Range("F35").Value = 41
$ws.Range("F36").Value = 58
$ws.Range("F38").Value = 96
$ws.Range("F39").Value = 3819
$ws.Range("F40").Value = 346
$ws.Range("F42").Value = 46
$ws.Range("F43").Value = 929
$ws.Range("F44").Value = 565
$ws.Range("F45").Value = 12
$ws.Range("F46").Value = 1471
$ws.Range("F48").Value = 569
$ws.Range("F49").Value = 660
$ws.Range("F50").Value = 12

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 46
$ws.Range("F17").Value = 230

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 7779
$ws.Range("F6").Value = 7779
$ws.Range("F7").Value = 7934
$ws.Range("F9").Value = 6794
$ws.Range("F10").Value = 3413
$ws.Range("F11").Value = 3750
$ws.Range("F12").Value = 56
$ws.Range("F13").Value = 51
$ws.Range("F14").Value = 77
$ws.Range("F15").Value = 85
$ws.Range("F16").Value = 481
$ws.Range("F17").Value = 46
$ws.Range("F18").Value = 63
$ws.Range("F21").Value = 339
$ws.Range("F22").Value = 3910
$ws.Range("F24").Value = 120
$ws.Range("F27").Value = 1013
$ws.Range("F28").Value = 343
$ws.Range("F29").Value = 1530
$ws.Range("F30").Value = 85
$ws.Range("F31").Value = 69
$ws.Range("F32").Value = 2809
$ws.Range("F33").Value = 1971
$ws.Range("F34").Value = 41
$ws.Range("F35").Value = 58
$ws.Range("F37").Value = 96
$ws.Range("F39").Value = 3819
$ws.Range("F40").Value = 346
$ws.Range("F42").Value = 46
$ws.Range("F43").Value = 929
$ws.Range("F44").Value = 565
$ws.Range("F45").Value = 230
$ws.Range("F46").Value = 1471
$ws.Range("F49").Value = 569
$ws.Range("F50").Value = 660

